# Apply updated numbers (想去人数 / 最低票价) to the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F3").Value = 1850
$ws1.Range("F4").Value = 483
$ws1.Range("F6").Value = 170
$ws1.Range("F7").Value = 2489
$ws1.Range("F8").Value = 155
$ws1.Range("F9").Value = 84
$ws1.Range("F10").Value = 170
$ws1.Range("F11").Value = 1505
$ws1.Range("F13").Value = 41
$ws1.Range("F14").Value = 324
$ws1.Range("G15").Value = 30
$ws1.Range("F16").Value = 21
$ws1.Range("F20").Value = 214
$ws1.Range("F21").Value = 7
$ws1.Range("F22").Value = 151
$ws1.Range("F24").Value = 1575
$ws1.Range("F26").Value = 389
$ws1.Range("F27").Value = 572
$ws1.Range("G27").Value = "已售罄"
$ws1.Range("F30").Value = 402

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F3").Value = 1850
$ws4.Range("F5").Value = 483
$ws4.Range("F7").Value = 170
$ws4.Range("F8").Value = 2489
$ws4.Range("F9").Value = 155
$ws4.Range("F10").Value = 84
$ws4.Range("F11").Value = 170
$ws4.Range("F12").Value = 1505
$ws4.Range("F14").Value = 41
$ws4.Range("F15").Value = 324
$ws4.Range("G16").Value = 30
$ws4.Range("F17").Value = 21
$ws4.Range("F21").Value = 214
$ws4.Range("F22").Value = 7
$ws4.Range("F23").Value = 151
$ws4.Range("F25").Value = 1575
$ws4.Range("F27").Value = 389
$ws4.Range("F28").Value = 572
$ws4.Range("G28").Value = "已售罄"
$ws4.Range("F31").Value = 402
